$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp title cell (A1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 13 de Junio de 2020 a las 21:24"

# Update province data rows (name + 4 numeric columns) to reflect the refreshed dataset
$ws.Cells.Item(4,1).Value = "Madrid"
$ws.Cells.Item(4,2).Value = 68852
$ws.Cells.Item(4,3).Value = 19425
$ws.Cells.Item(4,4).Value = 68852
$ws.Cells.Item(4,5).Value = 40736
$ws.Cells.Item(5,1).Value = "Cataluña"
$ws.Cells.Item(5,2).Value = 59019
$ws.Cells.Item(5,3).Value = 27229
$ws.Cells.Item(5,4).Value = 59019
$ws.Cells.Item(5,5).Value = 26203
$ws.Cells.Item(6,1).Value = "Castilla y Leon"
$ws.Cells.Item(6,2).Value = 18907
$ws.Cells.Item(6,3).Value = 8267
$ws.Cells.Item(6,4).Value = 18907
$ws.Cells.Item(6,5).Value = 8716
$ws.Cells.Item(7,1).Value = "Castilla-La Mancha"
$ws.Cells.Item(7,2).Value = 17259
$ws.Cells.Item(7,3).Value = 7922
$ws.Cells.Item(7,4).Value = 17259
$ws.Cells.Item(7,5).Value = 6392
$ws.Cells.Item(8,1).Value = "Pais Vasco"
$ws.Cells.Item(8,2).Value = 13156
$ws.Cells.Item(8,3).Value = 14646
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 1418
$ws.Cells.Item(9,1).Value = "Andalucia"
$ws.Cells.Item(9,2).Value = 12679
$ws.Cells.Item(9,3).Value = 604
$ws.Cells.Item(9,4).Value = 12679
$ws.Cells.Item(9,5).Value = 10671
$ws.Cells.Item(10,1).Value = "Bizkaia/Vizcaya"
$ws.Cells.Item(10,2).Value = 10332
$ws.Cells.Item(10,3).Value = 4423
$ws.Cells.Item(10,4).Value = 10332
$ws.Cells.Item(10,5).Value = 7124
$ws.Cells.Item(11,1).Value = "Galicia"
$ws.Cells.Item(11,2).Value = 9041
$ws.Cells.Item(11,3).Value = 8409
$ws.Cells.Item(11,4).Value = 28
$ws.Cells.Item(11,5).Value = 604
$ws.Cells.Item(12,1).Value = "Ciudad Real"
$ws.Cells.Item(12,2).Value = 6464
$ws.Cells.Item(12,3).Value = 3368
$ws.Cells.Item(12,4).Value = 6464
$ws.Cells.Item(12,5).Value = 2030
$ws.Cells.Item(13,1).Value = "Aragon"
$ws.Cells.Item(13,2).Value = 5695
$ws.Cells.Item(13,3).Value = 1097
$ws.Cells.Item(13,4).Value = 5695
$ws.Cells.Item(13,5).Value = 3772
$ws.Cells.Item(14,1).Value = "Valencia/Valencia"
$ws.Cells.Item(14,2).Value = 5609
$ws.Cells.Item(14,3).Value = 2767
$ws.Cells.Item(14,4).Value = 5609
$ws.Cells.Item(14,5).Value = 4907
$ws.Cells.Item(15,1).Value = "Zaragoza"
$ws.Cells.Item(15,2).Value = 5287
$ws.Cells.Item(15,3).Value = 1812
$ws.Cells.Item(15,4).Value = 5287
$ws.Cells.Item(15,5).Value = 2820
$ws.Cells.Item(16,1).Value = "Navarra"
$ws.Cells.Item(16,2).Value = 5259
$ws.Cells.Item(16,3).Value = 864
$ws.Cells.Item(16,4).Value = 5259
$ws.Cells.Item(16,5).Value = 3905
$ws.Cells.Item(17,1).Value = "Araba/Alava"
$ws.Cells.Item(17,2).Value = 4868
$ws.Cells.Item(17,3).Value = 4423
$ws.Cells.Item(17,4).Value = 4868
$ws.Cells.Item(17,5).Value = 7124
$ws.Cells.Item(18,1).Value = "Valladolid"
$ws.Cells.Item(18,2).Value = 4393
$ws.Cells.Item(18,3).Value = 2454
$ws.Cells.Item(18,4).Value = 4393
$ws.Cells.Item(18,5).Value = 1577
$ws.Cells.Item(19,1).Value = "Salamanca"
$ws.Cells.Item(19,2).Value = 4152
$ws.Cells.Item(19,3).Value = 2613
$ws.Cells.Item(19,4).Value = 4152
$ws.Cells.Item(19,5).Value = 1181
$ws.Cells.Item(20,1).Value = "La Rioja"
$ws.Cells.Item(20,2).Value = 4051
$ws.Cells.Item(20,3).Value = 583
$ws.Cells.Item(20,4).Value = 4051
$ws.Cells.Item(20,5).Value = 3107
$ws.Cells.Item(21,1).Value = "Toledo"
$ws.Cells.Item(21,2).Value = 3872
$ws.Cells.Item(21,3).Value = 1124
$ws.Cells.Item(21,4).Value = 3872
$ws.Cells.Item(21,5).Value = 1992
$ws.Cells.Item(22,1).Value = "Alacant/Alicante"
$ws.Cells.Item(22,2).Value = 3794
$ws.Cells.Item(22,3).Value = 1938
$ws.Cells.Item(22,4).Value = 3794
$ws.Cells.Item(22,5).Value = 3637
$ws.Cells.Item(23,1).Value = "Albacete"
$ws.Cells.Item(23,2).Value = 3775
$ws.Cells.Item(23,3).Value = 1872
$ws.Cells.Item(23,4).Value = 3775
$ws.Cells.Item(23,5).Value = 1389
$ws.Cells.Item(24,1).Value = "Leon"
$ws.Cells.Item(24,2).Value = 3569
$ws.Cells.Item(24,3).Value = 1533
$ws.Cells.Item(24,4).Value = 3569
$ws.Cells.Item(24,5).Value = 1629
$ws.Cells.Item(25,1).Value = "Segovia"
$ws.Cells.Item(25,2).Value = 3413
$ws.Cells.Item(25,3).Value = 2344
$ws.Cells.Item(25,4).Value = 3413
$ws.Cells.Item(25,5).Value = 868
$ws.Cells.Item(26,1).Value = "Gipuzkoa/Guipuzcoa"
$ws.Cells.Item(26,2).Value = 3116
$ws.Cells.Item(26,3).Value = 4423
$ws.Cells.Item(26,4).Value = 3116
$ws.Cells.Item(26,5).Value = 7124
$ws.Cells.Item(27,1).Value = "Extremadura"
$ws.Cells.Item(27,2).Value = 2919
$ws.Cells.Item(27,3).Value = 2422
$ws.Cells.Item(27,4).Value = 10
$ws.Cells.Item(27,5).Value = 487
$ws.Cells.Item(28,1).Value = "Malaga"
$ws.Cells.Item(28,2).Value = 2758
$ws.Cells.Item(28,3).Value = 185
$ws.Cells.Item(28,4).Value = 2758
$ws.Cells.Item(28,5).Value = 2295
$ws.Cells.Item(29,1).Value = "Burgos"
$ws.Cells.Item(29,2).Value = 2746
$ws.Cells.Item(29,3).Value = 1640
$ws.Cells.Item(29,4).Value = 2746
$ws.Cells.Item(29,5).Value = 901
$ws.Cells.Item(30,1).Value = "Asturias"
$ws.Cells.Item(30,2).Value = 2425
$ws.Cells.Item(30,3).Value = 1052
$ws.Cells.Item(30,4).Value = 2425
$ws.Cells.Item(30,5).Value = 1063
$ws.Cells.Item(31,1).Value = "Sevilla"
$ws.Cells.Item(31,2).Value = 2423
$ws.Cells.Item(31,3).Value = 308
$ws.Cells.Item(31,4).Value = 2423
$ws.Cells.Item(31,5).Value = 1840
$ws.Cells.Item(32,1).Value = "Granada"
$ws.Cells.Item(32,2).Value = 2413
$ws.Cells.Item(32,3).Value = 0
$ws.Cells.Item(32,4).Value = 2413
$ws.Cells.Item(32,5).Value = 2520
$ws.Cells.Item(33,1).Value = "Gran Canaria"
$ws.Cells.Item(33,2).Value = 2347
$ws.Cells.Item(33,3).Value = 659
$ws.Cells.Item(33,4).Value = 2347
$ws.Cells.Item(33,5).Value = 1537
$ws.Cells.Item(34,1).Value = "Soria"
$ws.Cells.Item(34,2).Value = 2290
$ws.Cells.Item(34,3).Value = 1774
$ws.Cells.Item(34,4).Value = 2290
$ws.Cells.Item(34,5).Value = 397
$ws.Cells.Item(35,1).Value = "Tenerife"
$ws.Cells.Item(35,2).Value = 2280
$ws.Cells.Item(35,3).Value = 623
$ws.Cells.Item(35,4).Value = 2280
$ws.Cells.Item(35,5).Value = 1506
$ws.Cells.Item(36,1).Value = "Cantabria"
$ws.Cells.Item(36,2).Value = 2246
$ws.Cells.Item(36,3).Value = 1981
$ws.Cells.Item(36,4).Value = 62
$ws.Cells.Item(36,5).Value = 203
$ws.Cells.Item(37,1).Value = "Caceres"
$ws.Cells.Item(37,2).Value = 1973
$ws.Cells.Item(37,3).Value = 66
$ws.Cells.Item(37,4).Value = 1973
$ws.Cells.Item(37,5).Value = 1505
$ws.Cells.Item(38,1).Value = "A Coruña"
$ws.Cells.Item(38,2).Value = 1969
$ws.Cells.Item(38,3).Value = 333
$ws.Cells.Item(38,4).Value = 1788
$ws.Cells.Item(38,5).Value = 67
$ws.Cells.Item(39,1).Value = "Avila"
$ws.Cells.Item(39,2).Value = 1935
$ws.Cells.Item(39,3).Value = 1179
$ws.Cells.Item(39,4).Value = 1935
$ws.Cells.Item(39,5).Value = 623
$ws.Cells.Item(40,1).Value = "Murcia"
$ws.Cells.Item(40,2).Value = 1587
$ws.Cells.Item(40,3).Value = 2180
$ws.Cells.Item(40,4).Value = 0
$ws.Cells.Item(40,5).Value = 148
$ws.Cells.Item(41,1).Value = "Pontevedra"
$ws.Cells.Item(41,2).Value = 1536
$ws.Cells.Item(41,3).Value = 333
$ws.Cells.Item(41,4).Value = 1411
$ws.Cells.Item(41,5).Value = 30
$ws.Cells.Item(42,1).Value = "Castello/Castellon"
$ws.Cells.Item(42,2).Value = 1486
$ws.Cells.Item(42,3).Value = 699
$ws.Cells.Item(42,4).Value = 1486
$ws.Cells.Item(42,5).Value = 1363
$ws.Cells.Item(43,1).Value = "Jaen"
$ws.Cells.Item(43,2).Value = 1387
$ws.Cells.Item(43,3).Value = 41
$ws.Cells.Item(43,4).Value = 1387
$ws.Cells.Item(43,5).Value = 1171
$ws.Cells.Item(44,1).Value = "Cordoba"
$ws.Cells.Item(44,2).Value = 1331
$ws.Cells.Item(44,3).Value = 0
$ws.Cells.Item(44,4).Value = 1331
$ws.Cells.Item(44,5).Value = 1350
$ws.Cells.Item(45,1).Value = "Guadalajara"
$ws.Cells.Item(45,2).Value = 1266
$ws.Cells.Item(45,3).Value = 644
$ws.Cells.Item(45,4).Value = 1266
$ws.Cells.Item(45,5).Value = 371
$ws.Cells.Item(46,1).Value = "Cuenca"
$ws.Cells.Item(46,2).Value = 1241
$ws.Cells.Item(46,3).Value = 339
$ws.Cells.Item(46,4).Value = 1241
$ws.Cells.Item(46,5).Value = 596
$ws.Cells.Item(47,1).Value = "Cadiz"
$ws.Cells.Item(47,2).Value = 1240
$ws.Cells.Item(47,3).Value = 535
$ws.Cells.Item(47,4).Value = 1240
$ws.Cells.Item(47,5).Value = 560
$ws.Cells.Item(48,1).Value = "Palencia"
$ws.Cells.Item(48,2).Value = 1205
$ws.Cells.Item(48,3).Value = 789
$ws.Cells.Item(48,4).Value = 1205
$ws.Cells.Item(48,5).Value = 333
$ws.Cells.Item(49,1).Value = "Huesca"
$ws.Cells.Item(49,2).Value = 1115
$ws.Cells.Item(49,3).Value = 544
$ws.Cells.Item(49,4).Value = 1115
$ws.Cells.Item(49,5).Value = 472
$ws.Cells.Item(50,1).Value = "Zamora"
$ws.Cells.Item(50,2).Value = 993
$ws.Cells.Item(50,3).Value = 586
$ws.Cells.Item(50,4).Value = 993
$ws.Cells.Item(50,5).Value = 322
$ws.Cells.Item(51,1).Value = "Badajoz"
$ws.Cells.Item(51,2).Value = 962
$ws.Cells.Item(51,3).Value = 0
$ws.Cells.Item(51,4).Value = 962
$ws.Cells.Item(51,5).Value = 1082
$ws.Cells.Item(52,1).Value = "Ourense"
$ws.Cells.Item(52,2).Value = 751
$ws.Cells.Item(52,3).Value = 333
$ws.Cells.Item(52,4).Value = 660
$ws.Cells.Item(52,5).Value = 22
$ws.Cells.Item(53,1).Value = "Teruel"
$ws.Cells.Item(53,2).Value = 664
$ws.Cells.Item(53,3).Value = 203
$ws.Cells.Item(53,4).Value = 664
$ws.Cells.Item(53,5).Value = 378
$ws.Cells.Item(54,1).Value = "Lugo"
$ws.Cells.Item(54,2).Value = 586
$ws.Cells.Item(54,3).Value = 333
$ws.Cells.Item(54,4).Value = 520
$ws.Cells.Item(54,5).Value = 11
$ws.Cells.Item(55,1).Value = "Almeria"
$ws.Cells.Item(55,2).Value = 498
$ws.Cells.Item(55,3).Value = 0
$ws.Cells.Item(55,4).Value = 498
$ws.Cells.Item(55,5).Value = 484
$ws.Cells.Item(56,1).Value = "Huelva"
$ws.Cells.Item(56,2).Value = 400
$ws.Cells.Item(56,3).Value = 0
$ws.Cells.Item(56,4).Value = 400
$ws.Cells.Item(56,5).Value = 391
$ws.Cells.Item(57,1).Value = "Mallorca"
$ws.Cells.Item(57,2).Value = 210
$ws.Cells.Item(57,3).Value = 18
$ws.Cells.Item(57,4).Value = 194
$ws.Cells.Item(57,5).Value = 12
$ws.Cells.Item(58,1).Value = "Ceuta"
$ws.Cells.Item(58,2).Value = 125
$ws.Cells.Item(58,3).Value = 98
$ws.Cells.Item(58,4).Value = 23
$ws.Cells.Item(58,5).Value = 4
$ws.Cells.Item(59,1).Value = "Melilla"
$ws.Cells.Item(59,2).Value = 121
$ws.Cells.Item(59,3).Value = 125
$ws.Cells.Item(59,4).Value = 0
$ws.Cells.Item(59,5).Value = 2
$ws.Cells.Item(60,1).Value = "La Palma"
$ws.Cells.Item(60,2).Value = 95
$ws.Cells.Item(60,3).Value = 22
$ws.Cells.Item(60,4).Value = 95
$ws.Cells.Item(60,5).Value = 68
$ws.Cells.Item(61,1).Value = "Lanzarote"
$ws.Cells.Item(61,2).Value = 84
$ws.Cells.Item(61,3).Value = 7
$ws.Cells.Item(61,4).Value = 84
$ws.Cells.Item(61,5).Value = 71
$ws.Cells.Item(62,1).Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Cells.Item(62,2).Value = 58
$ws.Cells.Item(62,3).Value = 0
$ws.Cells.Item(62,4).Value = 58
$ws.Cells.Item(62,5).Value = 3
$ws.Cells.Item(63,1).Value = "Fuerteventura"
$ws.Cells.Item(63,2).Value = 23
$ws.Cells.Item(63,3).Value = 0
$ws.Cells.Item(63,4).Value = 23
$ws.Cells.Item(63,5).Value = 42
$ws.Cells.Item(64,1).Value = "Ibiza"
$ws.Cells.Item(64,2).Value = 21
$ws.Cells.Item(64,3).Value = 18
$ws.Cells.Item(64,4).Value = 20
$ws.Cells.Item(64,5).Value = 1
$ws.Cells.Item(65,1).Value = "Menorca"
$ws.Cells.Item(65,2).Value = 15
$ws.Cells.Item(65,3).Value = 18
$ws.Cells.Item(65,4).Value = 13
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,1).Value = "La Gomera"
$ws.Cells.Item(66,2).Value = 8
$ws.Cells.Item(66,3).Value = 0
$ws.Cells.Item(66,4).Value = 8
$ws.Cells.Item(66,5).Value = 8
$ws.Cells.Item(67,1).Value = "Arroyo de la Luz"
$ws.Cells.Item(67,2).Value = 7
$ws.Cells.Item(67,3).Value = 0
$ws.Cells.Item(67,4).Value = 7
$ws.Cells.Item(67,5).Value = 0
$ws.Cells.Item(68,1).Value = "El Hierro"
$ws.Cells.Item(68,2).Value = 3
$ws.Cells.Item(68,3).Value = 0
$ws.Cells.Item(68,4).Value = 3
$ws.Cells.Item(68,5).Value = 3
$ws.Cells.Item(69,1).Value = "Formentera"
$ws.Cells.Item(69,2).Value = 0
$ws.Cells.Item(69,3).Value = 10
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 8
